# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The account-statement database (rows 16-46, sheet "Hoja1") is rebuilt:
#   - Column E ("Periodo Mora") is re-sorted from descending (2009 -> 1803)
#     to ascending (1803 -> 2009).
#   - Column F ("Valor Mora") gets the refreshed overdue amounts that go
#     along with "part 1" of the new account statement.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Ascending list of overdue periods (YYMM, stored as text) for E16:E46.
$periodos = @(
    "1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910",
    "1911","1912",
    "2001","2002","2003","2004","2005","2006","2007","2008","2009"
)

# Matching overdue values ("Valor Mora") for F16:F46.
$valores = @(
    20000,20000,20000,20000,20000,20000,20000,20000,20000,20000,
    30000,30000,30000,30000,30000,
    40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,
    37333
)

$firstRow = 16
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
